# Edit script: applies weekly price updates to the "Fruta, Femacal de
# La Calera - Caqui" sheet. Existing rows (2-35) get corrected Fecha
# (date), Calidad, Volumen, Precio minimo/maximo/promedio, Origen and
# Precio $/Kg values, and three brand-new rows (36-38) are appended
# at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---

# Row 2
$ws.Cells.Item(2,4).Value = 44315
$ws.Cells.Item(2,13).Value = 45

# Row 3
$ws.Cells.Item(3,4).Value = 44314
$ws.Cells.Item(3,12).Value = 'Primera'
$ws.Cells.Item(3,13).Value = 47

# Row 4
$ws.Cells.Item(4,4).Value = 44699
$ws.Cells.Item(4,12).Value = 'Especial'
$ws.Cells.Item(4,14).Value = 12000
$ws.Cells.Item(4,15).Value = 12000
$ws.Cells.Item(4,16).Value = 12000
$ws.Cells.Item(4,19).Value = 1200

# Row 5
$ws.Cells.Item(5,4).Value = 44699
$ws.Cells.Item(5,12).Value = 'Primera'
$ws.Cells.Item(5,13).Value = 60
$ws.Cells.Item(5,14).Value = 10000
$ws.Cells.Item(5,15).Value = 10000
$ws.Cells.Item(5,16).Value = 10000
$ws.Cells.Item(5,19).Value = 1000

# Row 6
$ws.Cells.Item(6,4).Value = 45091
$ws.Cells.Item(6,12).Value = 'Especial'
$ws.Cells.Item(6,13).Value = 54
$ws.Cells.Item(6,14).Value = 14000
$ws.Cells.Item(6,15).Value = 14000
$ws.Cells.Item(6,16).Value = 14000
$ws.Cells.Item(6,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(6,19).Value = 1400

# Row 7
$ws.Cells.Item(7,4).Value = 45091
$ws.Cells.Item(7,12).Value = 'Primera'
$ws.Cells.Item(7,13).Value = 58
$ws.Cells.Item(7,14).Value = 12000
$ws.Cells.Item(7,15).Value = 12000
$ws.Cells.Item(7,16).Value = 12000
$ws.Cells.Item(7,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(7,19).Value = 1200

# Row 8
$ws.Cells.Item(8,4).Value = 45091
$ws.Cells.Item(8,12).Value = 'Segunda'
$ws.Cells.Item(8,13).Value = 48
$ws.Cells.Item(8,18).Value = 'Región de O''Higgins'

# Row 9
$ws.Cells.Item(9,4).Value = 44302
$ws.Cells.Item(9,12).Value = 'Primera'
$ws.Cells.Item(9,13).Value = 45
$ws.Cells.Item(9,14).Value = 10000
$ws.Cells.Item(9,15).Value = 10000
$ws.Cells.Item(9,16).Value = 10000
$ws.Cells.Item(9,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(9,19).Value = 1000

# Row 10
$ws.Cells.Item(10,4).Value = 44323
$ws.Cells.Item(10,13).Value = 60
$ws.Cells.Item(10,14).Value = 10000
$ws.Cells.Item(10,15).Value = 10000
$ws.Cells.Item(10,16).Value = 10000
$ws.Cells.Item(10,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(10,19).Value = 1000

# Row 11
$ws.Cells.Item(11,4).Value = 44323
$ws.Cells.Item(11,13).Value = 50
$ws.Cells.Item(11,14).Value = 9000
$ws.Cells.Item(11,15).Value = 9000
$ws.Cells.Item(11,16).Value = 9000
$ws.Cells.Item(11,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(11,19).Value = 900

# Row 12
$ws.Cells.Item(12,4).Value = 44326
$ws.Cells.Item(12,12).Value = 'Primera'
$ws.Cells.Item(12,13).Value = 65
$ws.Cells.Item(12,14).Value = 10000
$ws.Cells.Item(12,15).Value = 10000
$ws.Cells.Item(12,16).Value = 10000
$ws.Cells.Item(12,19).Value = 1000

# Row 13
$ws.Cells.Item(13,4).Value = 44326
$ws.Cells.Item(13,12).Value = 'Segunda'
$ws.Cells.Item(13,13).Value = 67
$ws.Cells.Item(13,14).Value = 8000
$ws.Cells.Item(13,15).Value = 8000
$ws.Cells.Item(13,16).Value = 8000
$ws.Cells.Item(13,19).Value = 800

# Row 14
$ws.Cells.Item(14,4).Value = 44301

# Row 15
$ws.Cells.Item(15,4).Value = 44321
$ws.Cells.Item(15,13).Value = 58
$ws.Cells.Item(15,14).Value = 9000
$ws.Cells.Item(15,15).Value = 9000
$ws.Cells.Item(15,16).Value = 9000
$ws.Cells.Item(15,19).Value = 900

# Row 17
$ws.Cells.Item(17,4).Value = 44306

# Row 20
$ws.Cells.Item(20,4).Value = 44307
$ws.Cells.Item(20,13).Value = 40
$ws.Cells.Item(20,14).Value = 10000
$ws.Cells.Item(20,15).Value = 10000
$ws.Cells.Item(20,16).Value = 10000
$ws.Cells.Item(20,19).Value = 1000

# Row 21
$ws.Cells.Item(21,4).Value = 44322
$ws.Cells.Item(21,13).Value = 56

# Row 22
$ws.Cells.Item(22,4).Value = 44322
$ws.Cells.Item(22,13).Value = 40

# Row 23
$ws.Cells.Item(23,4).Value = 45082
$ws.Cells.Item(23,13).Value = 56
$ws.Cells.Item(23,14).Value = 15000
$ws.Cells.Item(23,15).Value = 15000
$ws.Cells.Item(23,16).Value = 15000
$ws.Cells.Item(23,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(23,19).Value = 1500

# Row 24
$ws.Cells.Item(24,4).Value = 45082
$ws.Cells.Item(24,13).Value = 67
$ws.Cells.Item(24,14).Value = 12000
$ws.Cells.Item(24,15).Value = 12000
$ws.Cells.Item(24,16).Value = 12000
$ws.Cells.Item(24,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(24,19).Value = 1200

# Row 25
$ws.Cells.Item(25,4).Value = 45082
$ws.Cells.Item(25,13).Value = 60
$ws.Cells.Item(25,14).Value = 10000
$ws.Cells.Item(25,15).Value = 10000
$ws.Cells.Item(25,16).Value = 10000
$ws.Cells.Item(25,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(25,19).Value = 1000

# Row 26
$ws.Cells.Item(26,4).Value = 44343
$ws.Cells.Item(26,12).Value = 'Especial'
$ws.Cells.Item(26,13).Value = 47
$ws.Cells.Item(26,14).Value = 10000
$ws.Cells.Item(26,15).Value = 10000
$ws.Cells.Item(26,16).Value = 10000
$ws.Cells.Item(26,18).Value = 'Región Metropolitana'
$ws.Cells.Item(26,19).Value = 1000

# Row 27
$ws.Cells.Item(27,4).Value = 44343
$ws.Cells.Item(27,12).Value = 'Primera'
$ws.Cells.Item(27,13).Value = 50
$ws.Cells.Item(27,14).Value = 9000
$ws.Cells.Item(27,15).Value = 9000
$ws.Cells.Item(27,16).Value = 9000
$ws.Cells.Item(27,18).Value = 'Región Metropolitana'
$ws.Cells.Item(27,19).Value = 900

# Row 28
$ws.Cells.Item(28,4).Value = 44343
$ws.Cells.Item(28,12).Value = 'Segunda'
$ws.Cells.Item(28,13).Value = 58
$ws.Cells.Item(28,14).Value = 8000
$ws.Cells.Item(28,15).Value = 8000
$ws.Cells.Item(28,16).Value = 8000
$ws.Cells.Item(28,18).Value = 'Región Metropolitana'
$ws.Cells.Item(28,19).Value = 800

# Row 29
$ws.Cells.Item(29,4).Value = 44308
$ws.Cells.Item(29,12).Value = 'Primera'
$ws.Cells.Item(29,13).Value = 45
$ws.Cells.Item(29,14).Value = 10000
$ws.Cells.Item(29,15).Value = 10000
$ws.Cells.Item(29,16).Value = 10000
$ws.Cells.Item(29,19).Value = 1000

# Row 30
$ws.Cells.Item(30,4).Value = 44308
$ws.Cells.Item(30,12).Value = 'Segunda'
$ws.Cells.Item(30,13).Value = 48
$ws.Cells.Item(30,14).Value = 8000
$ws.Cells.Item(30,15).Value = 8000
$ws.Cells.Item(30,16).Value = 8000
$ws.Cells.Item(30,19).Value = 800

# Row 31
$ws.Cells.Item(31,4).Value = 44319
$ws.Cells.Item(31,13).Value = 68

# Row 32
$ws.Cells.Item(32,4).Value = 44319
$ws.Cells.Item(32,12).Value = 'Segunda'
$ws.Cells.Item(32,13).Value = 57
$ws.Cells.Item(32,14).Value = 8000
$ws.Cells.Item(32,15).Value = 8000
$ws.Cells.Item(32,16).Value = 8000
$ws.Cells.Item(32,19).Value = 800

# Row 33
$ws.Cells.Item(33,4).Value = 44328
$ws.Cells.Item(33,12).Value = 'Primera'
$ws.Cells.Item(33,13).Value = 45

# Row 34
$ws.Cells.Item(34,4).Value = 44328
$ws.Cells.Item(34,12).Value = 'Segunda'
$ws.Cells.Item(34,13).Value = 48
$ws.Cells.Item(34,14).Value = 7000
$ws.Cells.Item(34,15).Value = 7000
$ws.Cells.Item(34,16).Value = 7000
$ws.Cells.Item(34,19).Value = 700

# Row 35
$ws.Cells.Item(35,4).Value = 44333
$ws.Cells.Item(35,12).Value = 'Especial'
$ws.Cells.Item(35,13).Value = 58
$ws.Cells.Item(35,14).Value = 10000
$ws.Cells.Item(35,15).Value = 10000
$ws.Cells.Item(35,16).Value = 10000
$ws.Cells.Item(35,19).Value = 1000

# --- Append new rows discovered in this weeks data pull ---

# Row 36
$ws.Cells.Item(36,1).Value = 3
$ws.Cells.Item(36,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(36,3).Value = 'Coquimbo'
$ws.Cells.Item(36,4).Value = 44333
$ws.Cells.Item(36,5).Value = 5
$ws.Cells.Item(36,6).Value = 'Fruta'
$ws.Cells.Item(36,7).Value = 100107
$ws.Cells.Item(36,8).Value = 'Otros'
$ws.Cells.Item(36,9).Value = 100107001
$ws.Cells.Item(36,10).Value = 'Caqui'
$ws.Cells.Item(36,11).Value = 'Mankaki'
$ws.Cells.Item(36,12).Value = 'Primera'
$ws.Cells.Item(36,13).Value = 65
$ws.Cells.Item(36,14).Value = 9000
$ws.Cells.Item(36,15).Value = 9000
$ws.Cells.Item(36,16).Value = 9000
$ws.Cells.Item(36,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(36,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(36,19).Value = 900
$ws.Cells.Item(36,20).Value = 10
$ws.Cells.Item(36,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 37
$ws.Cells.Item(37,1).Value = 3
$ws.Cells.Item(37,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(37,3).Value = 'Coquimbo'
$ws.Cells.Item(37,4).Value = 44333
$ws.Cells.Item(37,5).Value = 5
$ws.Cells.Item(37,6).Value = 'Fruta'
$ws.Cells.Item(37,7).Value = 100107
$ws.Cells.Item(37,8).Value = 'Otros'
$ws.Cells.Item(37,9).Value = 100107001
$ws.Cells.Item(37,10).Value = 'Caqui'
$ws.Cells.Item(37,11).Value = 'Mankaki'
$ws.Cells.Item(37,12).Value = 'Segunda'
$ws.Cells.Item(37,13).Value = 60
$ws.Cells.Item(37,14).Value = 8000
$ws.Cells.Item(37,15).Value = 8000
$ws.Cells.Item(37,16).Value = 8000
$ws.Cells.Item(37,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(37,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(37,19).Value = 800
$ws.Cells.Item(37,20).Value = 10
$ws.Cells.Item(37,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 38
$ws.Cells.Item(38,1).Value = 3
$ws.Cells.Item(38,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(38,3).Value = 'Coquimbo'
$ws.Cells.Item(38,4).Value = 44309
$ws.Cells.Item(38,5).Value = 5
$ws.Cells.Item(38,6).Value = 'Fruta'
$ws.Cells.Item(38,7).Value = 100107
$ws.Cells.Item(38,8).Value = 'Otros'
$ws.Cells.Item(38,9).Value = 100107001
$ws.Cells.Item(38,10).Value = 'Caqui'
$ws.Cells.Item(38,11).Value = 'Mankaki'
$ws.Cells.Item(38,12).Value = 'Primera'
$ws.Cells.Item(38,13).Value = 45
$ws.Cells.Item(38,14).Value = 10000
$ws.Cells.Item(38,15).Value = 10000
$ws.Cells.Item(38,16).Value = 10000
$ws.Cells.Item(38,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(38,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(38,19).Value = 1000
$ws.Cells.Item(38,20).Value = 10
$ws.Cells.Item(38,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

